# Updates cryptos list price (D) and volume (E) columns to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.577.08"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "3.441.08"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "3.439.51"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("D13").Value = "4.038.07"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "67.559.19"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000175"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "3.443.38"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("E31").Value = "  -4.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("E33").Value = "  -7.95%  "
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  -7.55%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.97%  "
$ws.Range("E43").Value = "  -5.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("D47").Value = "2.690.14"
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("E50").Value = "  -7.20%  "
$ws.Range("E51").Value = "  -5.31%  "
